$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.544.12"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").Value = "3.133.43"
$ws.Range("E3").Value = "  -3.39%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.75"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.62%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -7.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.68"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("E10").Value = "  -5.53%  "
$ws.Range("E11").Value = "  -3.44%  "
$ws.Range("D12").Value = "3.675.10"
$ws.Range("E12").Value = "  -3.85%  "
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "63.491.92"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.82"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.37%  "
$ws.Range("D16").Value = "3.128.00"
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("E17").Value = "  -3.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "405.04"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.03"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.31"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.199"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.479"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.61%  "
$ws.Range("E26").Value = "  -7.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.71"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.19%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E30").Value = "  -4.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.90"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.22"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.73"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.34%  "
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "152.66"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.31"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.20%  "
$ws.Range("D37").Value = "2.751.51"
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.64"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.35"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -8.92%  "
$ws.Range("E40").Value = "  -4.45%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.696"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.31%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0619"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.33"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.96%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0256"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "281.55"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -7.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "20.66"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.59%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0971"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.47"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.89"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -12.49%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.67"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.66%  "
